# Add the "InvalidLogin" worksheet (with sample invalid-login test data),
# placing it right after the existing "ValidLogin" sheet, and make it the
# active sheet/tab - mirroring a second data-source tab being added to the
# Selenium data workbook.

$wb = $excel.ActiveWorkbook
$validLogin = $wb.Worksheets.Item(1)

$invalidLogin = $wb.Worksheets.Add($null, $validLogin)
$invalidLogin.Name = "InvalidLogin"

# Populate column-by-column (A then B) so new shared-string entries are
# interned in header/row order: Usename, abcd, xyz, password.
$invalidLogin.Range("A1").Value = "Usename"
$invalidLogin.Range("A2").Value = "abcd"
$invalidLogin.Range("B2").Value = "xyz"
$invalidLogin.Range("B1").Value = "password"

# Leave the active cell on the username's password column, matching the
# saved selection state.
$invalidLogin.Range("B1").Select() | Out-Null
